$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Contraseña" (B) column first, then "Nombre de Usuario" (A),
# which previously only held the numeric row id (1..5) as a placeholder.
$ws.Range("B2").Value = "csgoisgoodgames123"
$ws.Range("B3").Value = "csgoisgoodgames124"
$ws.Range("B4").Value = "csgoisgoodgames125"
$ws.Range("B5").Value = "csgoisgoodgames126"
$ws.Range("B6").Value = "csgoisgoodgames127"

$ws.Range("A2").Value = "josegonzalezcoradopineed"
$ws.Range("A3").Value = "javiergonzalezcoradopineed"
$ws.Range("A4").Value = "davidgonzalezcoradopineed"
$ws.Range("A5").Value = "luisagonzalezcoradopineed"
$ws.Range("A6").Value = "silviagonzalezcoradopineed"

# Move the cursor/selection, matching where the author ended up after editing.
$ws.Range("E10").Select()
